# Rotate the "Recorded By" (column G) comma-separated list left by one
# position for every data row, e.g. "System, foo@bar.com" -> "foo@bar.com, System"
# Single-element values are left untouched (rotation of a 1-item list is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value2 = $rotated
    }
}
